# Weekly fruit/vegetable price update: insert a new weekly record as row 193
# (pushing the existing rows 193-216 down to 194-217, growing the used range
# from A1:R216 to A1:R217) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 193 - this shifts rows 193..216
# down to 194..217 and extends the sheet dimension automatically.
$ws.Rows.Item(193).Insert()

# Populate the newly-inserted row 193 with the new weekly price observation.
$ws.Range("A193").Value = 7
$ws.Range("B193").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C193").Value = "Ñuble"
$ws.Range("D193").Value = 45077
$ws.Range("E193").Value = 16
$ws.Range("F193").Value = 100112021
$ws.Range("G193").Value = "Ají"
$ws.Range("H193").Value = "Inferno"
$ws.Range("I193").Value = "Primera"
$ws.Range("J193").Value = 10
$ws.Range("K193").Value = 21000
$ws.Range("L193").Value = 21000
$ws.Range("M193").Value = 21000
$ws.Range("N193").Value = "`$/caja 15 kilos"
$ws.Range("O193").Value = "Región de Arica y Parinacota"
$ws.Range("P193").Value = 1400
$ws.Range("Q193").Value = 15
$ws.Range("R193").Value = "Hortaliza"
